{"js": "// Replace each two-digit multiplication expression in the document's\n// table cells with its new value, matched 1:1 in document order (per the\n// commit's regenerated practice sheet). Every \"old\" text below is unique\n// in the document, so an exact, case-sensitive search safely targets the\n// single run that holds it.\nconst replacements = [\n  [\"90\u00d717=\", \"68\u00d730=\"],\n  [\"58\u00d784=\", \"12\u00d745=\"],\n  [\"56\u00d779=\", \"88\u00d726=\"],\n  [\"59\u00d758=\", \"53\u00d767=\"],\n  [\"33\u00d785=\", \"65\u00d793=\"],\n  [\"75\u00d796=\", \"63\u00d732=\"],\n  [\"81\u00d727=\", \"62\u00d718=\"],\n  [\"22\u00d718=\", \"22\u00d786=\"],\n  [\"57\u00d791=\", \"41\u00d775=\"],\n  [\"36\u00d779=\", \"76\u00d772=\"],\n  [\"19\u00d750=\", \"53\u00d725=\"],\n  [\"27\u00d722=\", \"71\u00d785=\"],\n  [\"18\u00d790=\", \"67\u00d758=\"],\n  [\"63\u00d758=\", \"33\u00d755=\"],\n  [\"41\u00d712=\", \"92\u00d781=\"],\n  [\"96\u00d771=\", \"36\u00d763=\"],\n  [\"70\u00d766=\", \"96\u00d767=\"],\n  [\"99\u00d724=\", \"40\u00d738=\"],\n  [\"49\u00d795=\", \"12\u00d722=\"],\n  [\"19\u00d742=\", \"74\u00d741=\"],\n  [\"40\u00d755=\", \"89\u00d719=\"],\n  [\"44\u00d716=\", \"47\u00d755=\"],\n  [\"25\u00d769=\", \"44\u00d725=\"],\n  [\"20\u00d782=\", \"90\u00d794=\"],\n  [\"36\u00d799=\", \"21\u00d747=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find expression \"${oldText}\" to replace.`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit multiplication expression in the document's\n# table cells with its new value, matched 1:1 in document order (per the\n# commit's regenerated practice sheet). Every \"old\" text below is unique\n# in the document, so Find/Replace safely targets exactly one run each.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"90\u00d717=\", \"68\u00d730=\"),\n    @(\"58\u00d784=\", \"12\u00d745=\"),\n    @(\"56\u00d779=\", \"88\u00d726=\"),\n    @(\"59\u00d758=\", \"53\u00d767=\"),\n    @(\"33\u00d785=\", \"65\u00d793=\"),\n    @(\"75\u00d796=\", \"63\u00d732=\"),\n    @(\"81\u00d727=\", \"62\u00d718=\"),\n    @(\"22\u00d718=\", \"22\u00d786=\"),\n    @(\"57\u00d791=\", \"41\u00d775=\"),\n    @(\"36\u00d779=\", \"76\u00d772=\"),\n    @(\"19\u00d750=\", \"53\u00d725=\"),\n    @(\"27\u00d722=\", \"71\u00d785=\"),\n    @(\"18\u00d790=\", \"67\u00d758=\"),\n    @(\"63\u00d758=\", \"33\u00d755=\"),\n    @(\"41\u00d712=\", \"92\u00d781=\"),\n    @(\"96\u00d771=\", \"36\u00d763=\"),\n    @(\"70\u00d766=\", \"96\u00d767=\"),\n    @(\"99\u00d724=\", \"40\u00d738=\"),\n    @(\"49\u00d795=\", \"12\u00d722=\"),\n    @(\"19\u00d742=\", \"74\u00d741=\"),\n    @(\"40\u00d755=\", \"89\u00d719=\"),\n    @(\"44\u00d716=\", \"47\u00d755=\"),\n    @(\"25\u00d769=\", \"44\u00d725=\"),\n    @(\"20\u00d782=\", \"90\u00d794=\"),\n    @(\"36\u00d799=\", \"21\u00d747=\")\n)\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $oldText,    # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        $wdFindContinue,  # Wrap\n        $false,      # Format\n        $newText,    # ReplaceWith\n        $wdReplaceAll     # Replace\n    ) | Out-Null\n}\n"}
